$d = $word.ActiveDocument

# Locate the "CAUSAS QUE HAN PROVOCADO EL ACCIDENTE:" paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("CAUSAS QUE HAN PROVOCADO EL ACCIDENTE:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'CAUSAS QUE HAN PROVOCADO EL ACCIDENTE:' paragraph"
}
$headingStart = $rng.Start

$startIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $headingStart) {
        $startIdx = $i
        break
    }
}
if ($startIdx -eq -1) {
    throw "Could not resolve the heading paragraph index"
}

# Delete the heading paragraph and the blank paragraph immediately after it
# (both entirely, including their paragraph marks), so whatever was before
# them now runs straight into the "{lista_causas}" paragraph.
$startPara = $d.Paragraphs.Item($startIdx)
$endPara = $d.Paragraphs.Item($startIdx + 1)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# Remove the "lista_causas" bookmark together with the literal
# "{lista_causas}" placeholder text it wraps, but keep the paragraph (and
# the page-break run that follows it) intact.
$bm = $d.Bookmarks.Item("lista_causas")
$bmRange = $d.Range($bm.Start, $bm.End)
$bmRange.Text = ""
$bm.Delete()
